# v1.14: Solapa Todas las transacciones, modal edición completa con combos
$wb = $excel.ActiveWorkbook

# --- Sheet "Log": append new bitacora entry as row 45 ---
$wsLog = $wb.Worksheets.Item("Log")
$wsLog.Cells.Item(45, 1).Value = "28/02/2025"
$wsLog.Cells.Item(45, 2).Value = "12:00"
$wsLog.Cells.Item(45, 3).Value = "Solapa Todas las transacciones y edición completa"
$wsLog.Cells.Item(45, 4).Value = "Nueva solapa Todas las transacciones: listado con todas las columnas, filtros por mes y categoría, botón Editar por registro. Modal de edición ampliado: todos los campos editables; combos para valores normalizados (categoría, cuenta contable, tipo movimiento, status, medio pago, moneda, origen archivo). editado y editado_detalle al guardar."
$wsLog.Cells.Item(45, 5).Value = "Implementacion"

# --- Sheet "Resumen": append two new functionality rows 41-42 ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Cells.Item(41, 1).Value = "Todas las transacciones"
$wsResumen.Cells.Item(41, 2).Value = "Solapa que lista todas las transacciones con todas las columnas. Filtros por mes y categoría. Botón Editar por registro abre modal de edición completa."
$wsResumen.Cells.Item(42, 1).Value = "Edición completa de registros"
$wsResumen.Cells.Item(42, 2).Value = "Modal de edición con todos los campos: fecha, mes, año, tipo movimiento, monto, moneda, status, medio pago, categoría, cuenta contable, origen archivo, descripción, cliente, cat_desc, id_origen, id_operación. Combos para campos normalizados (valores existentes en BD). editado y editado_detalle al guardar."

# --- Sheet "Versiones": append new version row 16 ---
# Leading apostrophe keeps "1.14" stored as text (matches the existing
# version numbers in column A, which are text, not numbers).
$wsVersiones = $wb.Worksheets.Item("Versiones")
$wsVersiones.Cells.Item(16, 1).Value = "'1.14"
$wsVersiones.Cells.Item(16, 2).Value = "28/02/2025"
$wsVersiones.Cells.Item(16, 3).Value = "Solapa Todas las transacciones (filtros mes y categoría); modal edición completa con todos los campos y combos para normalizados"

# --- Sheet "Presupuesto": insert new budget line at row 7, shifting existing rows down ---
$wsPresupuesto = $wb.Worksheets.Item("Presupuesto")
$wsPresupuesto.Rows.Item(7).Insert()
$wsPresupuesto.Cells.Item(7, 1).Value = "Listado y edición completa de transacciones"
$wsPresupuesto.Cells.Item(7, 2).Value = "Solapa Todas las transacciones con listado completo, filtros por mes y categoría, y modal de edición con todos los campos y combos para valores normalizados (categoría, cuenta contable, tipo movimiento, status, medio pago, moneda, origen archivo)."
$wsPresupuesto.Cells.Item(7, 3).Value = 45000
